$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose literal text would otherwise be auto-detected as a number by
# Excel (losing the exact original text, e.g. trailing zeros). The source
# workbook stores every Price/Volume cell as text, so force Text format on
# these before assigning the value to preserve the literal string.
$ws.Range('D2').Value = '42.684.12'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '2.304.62'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.69'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.61'
$ws.Range('E6').Value = '  -5.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -4.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.504'
$ws.Range('E9').Value = '  -4.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.72'
$ws.Range('E10').Value = '  -5.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').Value = '  -2.92%  '
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.77'
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').Value = '2.662.44'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.71'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').Value = '2.334.60'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.805'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '42.605.91'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').Value = '0.0₃0906'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.53'
$ws.Range('E20').Value = '  -4.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.06'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.09'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.38'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  -4.64%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.06'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  +6.72%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.78'
$ws.Range('E29').Value = '  -5.44%  '
$ws.Range('E30').Value = '  -4.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.57'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.03'
$ws.Range('E33').Value = '  -4.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.58'
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('E35').Value = '  -5.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.98'
$ws.Range('E36').Value = '  -7.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0707'
$ws.Range('E37').Value = '  -4.51%  '
$ws.Range('E38').Value = '  -5.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.81'
$ws.Range('E39').Value = '  -4.27%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  -5.25%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.111'
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.54'
$ws.Range('E42').Value = '  -1.77%  '
$ws.Range('D43').Value = '1.971.80'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0281'
$ws.Range('E44').Value = '  -3.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.45'
$ws.Range('E45').Value = '  -5.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.25'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('E47').Value = '  -6.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.34'
$ws.Range('E48').Value = '  -4.59%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.529.10'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.82'
$ws.Range('E50').Value = '  -2.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.69'
$ws.Range('E51').Value = '  -0.36%  '
